# The commit swaps the contents of ppt/theme/theme1.xml (the "Office
# Theme" default palette) and ppt/theme/theme2.xml (the "Integral" /
# "Red Violet" palette used by the slide master) — i.e. the deck's
# live theme (theme2.xml, the one driving every slide through
# SlideMaster) switches from the pink/purple "Integral" palette to the
# plain "Office Theme" colors.
#
# The PowerPoint object model exposes exactly one editable theme
# resource from a deck with a single slide master — reached via
# SlideMaster.Theme.ThemeColorScheme (NotesMaster.Theme resolves to
# the same 12-slot color table). We drive that theme's 12 color slots
# (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink) to the "Office
# Theme" values that theme1.xml held before the edit, matching the
# diff's end state for the deck's active theme.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$colors = $master.Theme.ThemeColorScheme

# Office Theme palette (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink),
# expressed as the VBA RGB() integer encoding (R + G*256 + B*65536)
# that PowerPoint's ThemeColor.RGB property uses.
$colors.Item(1).RGB  = 0x000000   # dk1      000000
$colors.Item(2).RGB  = 0xFFFFFF   # lt1      FFFFFF
$colors.Item(3).RGB  = 0x6A5444   # dk2      44546A
$colors.Item(4).RGB  = 0xE6E6E7   # lt2      E7E6E6
$colors.Item(5).RGB  = 0xD59B5B   # accent1  5B9BD5
$colors.Item(6).RGB  = 0x317DED   # accent2  ED7D31
$colors.Item(7).RGB  = 0xA5A5A5   # accent3  A5A5A5
$colors.Item(8).RGB  = 0x00C0FF   # accent4  FFC000
$colors.Item(9).RGB  = 0xC47244   # accent5  4472C4
$colors.Item(10).RGB = 0x47AD70   # accent6  70AD47
$colors.Item(11).RGB = 0xC16305   # hlink    0563C1
$colors.Item(12).RGB = 0x724F95   # folHlink 954F72
